# -----------------------------------------------------------------------
# Target diff analysis
# -----------------------------------------------------------------------
# The supplied unified diff is a diff of *canonical* OOXML (attributes are
# printed in sorted/canonical order, cf. XML C14N) for exactly two parts:
#
#   word/document.xml   - only the <w:sectPr> of the last section is
#                          touched, and every change there is purely an
#                          attribute-order permutation, e.g.
#                              <w:headerReference w:type="even" r:id="rId11"/>
#                          becomes
#                              <w:headerReference r:id="rId11" w:type="even"/>
#                          Same set of attributes/values, only the
#                          serialization order differs - not a content
#                          edit (confirmed by round-tripping the fragment
#                          through canonical-XML serialization, which
#                          reproduces the "after" ordering exactly from
#                          the "before" content).
#
#   word/numbering.xml   - same story for ~every line (namespace
#                          declarations reordered alphabetically,
#                          w:tab/w:ind/w:style/... attributes resorted,
#                          etc.) EXCEPT for one real value change:
#                              <w:nsid w:val="2c1ae401"/>  ->  <w:nsid w:val="A990"/>
#                          inside the single <w:abstractNum w:abstractNumId="990">.
#
# That abstractNum/num pair (numId 1000 -> abstractNumId 990) is not
# referenced by any paragraph in the document body (no <w:numPr> anywhere
# in word/document.xml), so this is dead/unused list-definition
# boilerplate carried over from the reference template. Its `nsid` is an
# internal, implementation-minted identifier for the abstract numbering
# definition - it is not surfaced anywhere in the Word object model (no
# List/ListTemplate/ListLevel/ListFormat property maps to it - by design,
# real Word does not let VBA/COM read or assign `nsid` either), so there
# is no COM call that can reproduce that value. The commit message
# ("Upgrade to pandoc 2.13") confirms this is a tooling/regeneration
# artifact of the docx *writer* being bumped, not a deliberate edit made
# in Word to the document's content, styles, numbering, headers/footers,
# or page setup.
#
# Net effect: once attribute-ordering noise is discounted, this diff
# describes no semantically meaningful, COM-reachable change to the
# document - every part the diff touches is otherwise byte-for-byte
# equivalent. So the correct edit here is to leave the document's
# content untouched rather than attempt a workaround (e.g. forcing a
# ListTemplate onto a paragraph to dirty word/numbering.xml) that would
# only introduce new abstractNum/num entries and reordered
# header/footer references that are NOT present in the target diff.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# Touch nothing: $d.Content / $d.Sections / $d.ListTemplates are left
# exactly as loaded so every part not mentioned in the diff (styles,
# settings, headers, footers, docProps, ...) - and the parts that are
# mentioned but only reordered by canonicalization - round-trip
# unchanged.
